$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of table-name metadata below the existing data (row 76)
$ws.Cells.Item(76, 1).Value  = "table_name"
$ws.Cells.Item(76, 4).Value  = "PVLT_Can_Han_Canadian_Hanwha"
$ws.Cells.Item(76, 2).Value  = "dbo.PVLT_Pan_LG_Panasonic_LG"
$ws.Cells.Item(76, 3).Value  = "dbo.PVLT_Pan_LG_Panasonic_LG"
$ws.Cells.Item(76, 5).Value  = "PVLT_Can_Han_Canadian_Hanwha"
$ws.Cells.Item(76, 7).Value  = "dbo.PVLT_TrinaJinko_Jinko"
$ws.Cells.Item(76, 8).Value  = "dbo.PVLT_MissionSolar_MissionSolar"
$ws.Cells.Item(76, 9).Value  = "PVLT_Can_Han_Canadian_Hanwha"
$ws.Cells.Item(76, 10).Value = "PVLT_Can_Han_Canadian_Hanwha"

# Update the saved scroll position / selection of the sheet view
$win = $wb.Windows.Item(1)
$ws.Range("D56").Select()
$win.ScrollRow = 2
$win.ScrollColumn = 1
